$wb = $excel.ActiveWorkbook

# --- Insert the new "2022-Q1" sheet between "2021-Q4" and "总计" ---
$q4 = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Header row (bold, bordered, centered - matches the other data sheet's header style)
$header = $ws.Range("B1:H1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data row 2 - first column is the (styled) numeric row index
$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Value = 0

# B2:G2 are stored as plain text (force text type, then drop the format noise)
$ws.Range("B2:G2").NumberFormat = "@"
$ws.Range("B2").Value = "001743"
$ws.Range("C2").Value = "诺安优选回报灵活配置混合"
$ws.Range("D2").Value = "6.13"
$ws.Range("E2").Value = "71.32"
$ws.Range("F2").Value = "5.63"
$ws.Range("G2").Value = "0.3451"
$ws.Range("B2:G2").ClearFormats()

# H2 is a genuine number (rank)
$ws.Range("H2").Value = 4

# --- Update the "总计" (totals) sheet: insert a new top row for 2022-Q1 ---
$tot = $wb.Worksheets.Item("总计")
$tot.Range("A2:D2").Insert()
$tot.Range("B2:D2").ClearFormats()

$a2b = $tot.Range("A2")
$a2b.Font.Bold = $true
$a2b.Borders.LineStyle = 1
$a2b.HorizontalAlignment = -4108
$a2b.VerticalAlignment = -4160
$a2b.Value = 0

$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 1
$tot.Range("D2").Value = 0.35

# Renumber the old first data row's index cell (now row 3) back to 1
$a3 = $tot.Range("A3")
$a3.Font.Bold = $true
$a3.Borders.LineStyle = 1
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4160
$a3.Value = 1

# Restore the originally active sheet/selection
$q4.Activate() | Out-Null
$q4.Range("A1").Select() | Out-Null
